$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is an "Add Fields" helper listing example field values under a
# NOMBRE header. Relabel the header's sample value (was "ALGODON") to
# "TINTE" and clear out the other example rows (NYLON/SEDA/POLIESTER),
# leaving a blank template list A3:A6 that all share the same formatting.
$ws.Range("A2").Value = "TINTE"
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()

# Re-apply the shared look (Arial 11, black, centered) across A3:A6 so the
# now-empty rows are uniformly formatted (drops the stray duplicate style).
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 11
$ws.Range("A4").Font.Color = 0
$ws.Range("A6").Font.Name = "Arial"
$ws.Range("A6").Font.Size = 11
$ws.Range("A6").Font.Color = 0
$ws.Range("A3").Font.Name = "Arial"
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Color = 0
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 11
$ws.Range("A5").Font.Color = 0

# Tiny incidental default-width recalculation that tagged along with the
# font cleanup above (best effort; harmless if the host ignores it).
$ws.StandardWidth = 12.71484375

# Move the selection up onto the now-blank entry rows.
$ws.Range("A3:A5").Select() | Out-Null
